# scale data push and pull
#
# Rescale the "value" column (D) by multiplying every numeric value in
# D2:D33 by 10000 (the source unit changed from 10^4 CNY ("wan yuan") to
# plain CNY). Empty cells (e.g. D29, which has no observation) are left
# untouched.
#
# NOTE: this runtime's PowerShell functions do not have their own local
# variable scope (callee writes are visible to the caller, much like
# dot-sourcing), so every variable name below is kept unique across the
# whole script/functions to avoid accidental aliasing with the outer loop.
#
# The rescale itself is performed as an exact decimal-point shift of the
# original textual value (equivalent to multiplying the exact decimal
# value by 10000) rather than an IEEE-754 double multiplication, so the
# resulting stored value matches bit-for-bit what a precise decimal-aware
# tool would produce (plain double multiplication can be off by 1 ULP).

function ScaleValueText([string]$txt) {
    $zerosTable = @("", "0", "00", "000", "0000")

    $isNeg = $false
    if ($txt.StartsWith("-")) {
        $isNeg = $true
        $txt = $txt.Substring(1)
    }

    if ($txt.Contains(".")) {
        $txtParts = $txt -split '\.', 2
        $txtInt = $txtParts[0]
        $txtFrac = $txtParts[1]
    } else {
        $txtInt = $txt
        $txtFrac = ""
    }

    # Move the decimal point 4 places to the right (multiply by 10000
    # exactly, at the text level).
    $takeLen = [Math]::Min(4, $txtFrac.Length)
    $movedDigits = $txtFrac.Substring(0, $takeLen)
    $leftoverFrac = $txtFrac.Substring($takeLen)
    $padLen = 4 - $takeLen
    $padZeros = $zerosTable[$padLen]

    $combinedInt = "$txtInt$movedDigits$padZeros"
    $combinedInt = $combinedInt.TrimStart("0")
    if ($combinedInt -eq "") {
        $combinedInt = "0"
    }

    if ($leftoverFrac -ne "") {
        $finalText = "$combinedInt.$leftoverFrac"
    } else {
        $finalText = $combinedInt
    }

    if ($isNeg) {
        $finalText = "-$finalText"
    }

    return $finalText
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $targetCell = $ws.Cells.Item($row, 4)
    $cellVal = $targetCell.Value2
    if ($cellVal -ne $null -and $cellVal -ne "") {
        $cellValText = $cellVal.ToString("R")
        $scaledText = ScaleValueText $cellValText
        $targetCell.Value2 = [double]$scaledText
    }
}
